# Apply the "add experimental" / date-bump edit described by the commit:
#   - CodeSystem.experimental is now set -> Metadata sheet row "Experimental"
#     gets a Value of literal text "true" (was previously blank).
#   - The publication Date moves from 2025-01-22T15:00:55+00:00
#     to 2025-01-28T15:58:19+00:00.
#
# The Metadata sheet is laid out as two columns (A: Property, B: Value),
# one property per row; "Experimental" is row 7 and "Date" is row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: Experimental -> true ---------------------------------------
# Writing the bare text "true"/"'true" to a cell makes Excel (and this
# engine, faithfully) auto-coerce it to a Boolean TRUE, which would store
# the cell as t="b" instead of a plain shared-string "true" and would also
# mint a brand-new quote-prefixed cell style. Routing it through a text
# formula first, then collapsing the formula down to its cached value via
# a Copy / PasteSpecial(values-only), keeps the literal text "true" while
# preserving the existing cell style untouched.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues

# --- Row 8: Date -> new timestamp --------------------------------------
$ws.Range("B8").Value = "2025-01-28T15:58:19+00:00"
